$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 21.84976866666667
$ws.Cells.Item(2, 8).Value = 65.549306
$ws.Cells.Item(2, 9).Value = 0.05020018890879543
$ws.Cells.Item(2, 10).Value = 0.05020018890879543
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.095195666666667
$ws.Cells.Item(2, 14).Value = 6.285587
$ws.Cells.Item(2, 15).Value = 0.8546922300706357
$ws.Cells.Item(2, 16).Value = 0.8546922300706358
$ws.Cells.Item(2, 17).Value = 45.77954062806911
$ws.Cells.Item(2, 18).Value = 412.015865652622
$ws.Cells.Item(2, 19).Value = 0.04290571140842556
$ws.Cells.Item(2, 20).Value = 0.04290571140842556

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 21.84976866666667
$ws.Cells.Item(3, 8).Value = 65.549306
$ws.Cells.Item(3, 9).Value = 0.05020018890879543
$ws.Cells.Item(3, 10).Value = 0.05020018890879543
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.356208
$ws.Cells.Item(3, 14).Value = 1.068624
$ws.Cells.Item(3, 15).Value = 0.1453077699293643
$ws.Cells.Item(3, 16).Value = 0.1453077699293643
$ws.Cells.Item(3, 17).Value = 7.783062397216001
$ws.Cells.Item(3, 18).Value = 70.047561574944
$ws.Cells.Item(3, 19).Value = 0.00729447750036987
$ws.Cells.Item(3, 20).Value = 0.00729447750036987

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 385.0524703333334
$ws.Cells.Item(4, 8).Value = 1155.157411
$ws.Cells.Item(4, 9).Value = 0.8846641374295412
$ws.Cells.Item(4, 10).Value = 0.8846641374295412
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.095195666666667
$ws.Cells.Item(4, 14).Value = 6.285587
$ws.Cells.Item(4, 15).Value = 0.8546922300706357
$ws.Cells.Item(4, 16).Value = 0.8546922300706358
$ws.Cells.Item(4, 17).Value = 806.7602672816953
$ws.Cells.Item(4, 18).Value = 7260.842405535257
$ws.Cells.Item(4, 19).Value = 0.7561155644831699
$ws.Cells.Item(4, 20).Value = 0.75611556448317

# Row 5 (new)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fn1"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 385.0524703333334
$ws.Cells.Item(5, 8).Value = 1155.157411
$ws.Cells.Item(5, 9).Value = 0.8846641374295412
$ws.Cells.Item(5, 10).Value = 0.8846641374295412
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.356208
$ws.Cells.Item(5, 14).Value = 1.068624
$ws.Cells.Item(5, 15).Value = 0.1453077699293643
$ws.Cells.Item(5, 16).Value = 0.1453077699293643
$ws.Cells.Item(5, 17).Value = 137.158770352496
$ws.Cells.Item(5, 18).Value = 1234.428933172464
$ws.Cells.Item(5, 19).Value = 0.1285485729463713
$ws.Cells.Item(5, 20).Value = 0.1285485729463713

# Row 6 (new)
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Fn1"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 28.350479
$ws.Cells.Item(6, 8).Value = 85.05143699999999
$ws.Cells.Item(6, 9).Value = 0.06513567366166337
$ws.Cells.Item(6, 10).Value = 0.06513567366166337
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.095195666666667
$ws.Cells.Item(6, 14).Value = 6.285587
$ws.Cells.Item(6, 15).Value = 0.8546922300706357
$ws.Cells.Item(6, 16).Value = 0.8546922300706358
$ws.Cells.Item(6, 17).Value = 59.39980074872432
$ws.Cells.Item(6, 18).Value = 534.5982067385189
$ws.Cells.Item(6, 19).Value = 0.05567095417904023
$ws.Cells.Item(6, 20).Value = 0.05567095417904024

# Row 7 (new)
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Fn1"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 28.350479
$ws.Cells.Item(7, 8).Value = 85.05143699999999
$ws.Cells.Item(7, 9).Value = 0.06513567366166337
$ws.Cells.Item(7, 10).Value = 0.06513567366166337
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.356208
$ws.Cells.Item(7, 14).Value = 1.068624
$ws.Cells.Item(7, 15).Value = 0.1453077699293643
$ws.Cells.Item(7, 16).Value = 0.1453077699293643
$ws.Cells.Item(7, 17).Value = 10.098667423632
$ws.Cells.Item(7, 18).Value = 90.88800681268799
$ws.Cells.Item(7, 19).Value = 0.009464719482623134
$ws.Cells.Item(7, 20).Value = 0.009464719482623134
